$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Moods")

for ($row = 2; $row -le 111; $row++) {
    $cell = $ws.Cells.Item($row, 8)  # column H
    $cell.Formula = '="UPDATE MOODS SET MorningWeight="&C' + $row + '&", AfternoonWeight="&D' + $row + '&", EveningWeight="&E' + $row + '&", LateWeight="&F' + $row + '&" WHERE Name=' + "'" + '"&B' + $row + '&"' + "'" + ';"'
}

# Apply the style used for s="4" (general alignment explicitly applied) to H2:H111
$range = $ws.Range("H2:H111")
$range.HorizontalAlignment = 1  # xlHAlignGeneral, forces explicit alignment attrs similar to s=4

# Update frozen pane scroll position and selection on the active sheet view
$ws.Activate()
$appWin = $excel.ActiveWindow
$appWin.ScrollRow = 103
$ws.Range("H2:H111").Select()
